# Add a new "Exceptional items" column to the Quarterly sheet.
#
# The column is inserted between "P/l before exceptional items & tax" (K)
# and "P/l before tax" (old L, now shifted to M), pushing every column from
# the old L through T one place to the right (new M through U).
#
# Two quarters (Dec '13 -> row 3, and Mar '15 -> row 11) have an actual
# "Exceptional items" figure; every other quarter has none.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a new column before column L; this shifts the old L:T into M:U
# and carries along the formatting of the neighbouring column.
$ws.Columns("L:L").Insert()

# Header rows for the new column.
$ws.Range("L1").Value = "Exceptional items"
$ws.Range("L2").Value = "Exceptional Items"

# The only two quarters that actually report an exceptional item.
$ws.Range("L3").Value = 120
$ws.Range("L11").Value = 61.27
